# Update to "Australia ALeague" odds tables
# 1) Rows 112 and 113 (match ids 110/111, B column = external "id" field) have
#    their entire row content (columns B..AC) swapped with each other. Column A
#    (the local/sequential id) stays where it is.
# 2) A handful of odds cells in rows 146, 147, 149, 150 and 151 are refreshed
#    to newer quoted values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: swap row 112 <-> row 113 across columns B (2) .. AC (29)
# ---------------------------------------------------------------------------
$firstCol = 2   # column B
$lastCol  = 29  # column AC
$rowA = 112
$rowB = 113

for ($c = $firstCol; $c -le $lastCol; $c++) {
    $cellA = $ws.Cells.Item($rowA, $c)
    $cellB = $ws.Cells.Item($rowB, $c)

    $valA = $cellA.Value()
    $valB = $cellB.Value()

    $cellA.Value = $valB
    $cellB.Value = $valA
}

# ---------------------------------------------------------------------------
# Part 2: refresh individual odds values
# ---------------------------------------------------------------------------

# Row 146
$ws.Range("N146").Value = 3.5
$ws.Range("P146").Value = 2.05
$ws.Range("R146").Value = 2.08
$ws.Range("S146").Value = 1.82

# Row 147
$ws.Range("N147").Value = 1.8
$ws.Range("P147").Value = 4
$ws.Range("Q147").Value = -0.75
$ws.Range("R147").Value = 2.03
$ws.Range("S147").Value = 1.87
$ws.Range("U147").Value = 1.875
$ws.Range("V147").Value = 1.975

# Row 149
$ws.Range("N149").Value = 1.8
$ws.Range("O149").Value = 4.2
$ws.Range("P149").Value = 3.8
$ws.Range("Q149").Value = -0.75
$ws.Range("R149").Value = 2.04
$ws.Range("S149").Value = 1.86
$ws.Range("U149").Value = 1.9
$ws.Range("V149").Value = 1.95

# Row 150
$ws.Range("R150").Value = 1.84
$ws.Range("S150").Value = 2.06

# Row 151
$ws.Range("N151").Value = 1.4
$ws.Range("P151").Value = 6.5
$ws.Range("R151").Value = 1.87
$ws.Range("S151").Value = 2.03
$ws.Range("U151").Value = 2
$ws.Range("V151").Value = 1.85
